# Generate Report for Handoff
# Refreshes the localization-status report: the overall status moves from
# "In Translation" to "Ready for handoff" and the "latest handoff" / "latest
# HO xliff generate" timestamps are bumped to the new run time. The Status /
# Datetime columns also widen slightly to fit the new text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" ------------------------
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- Timestamps bumped forward by the new handoff run -----------------------
# Overview!G2 "Latest HO Xliff Generate Date" and de-de!H2 "Latest Handoff
# Datetime" shared the same original value.
$overview.Range("G2").Value = "2016-09-05 16:46:45"
$dede.Range("H2").Value     = "2016-09-05 16:46:45"

# zh-cn!H2 "Latest Handoff Datetime"
$zhcn.Range("H2").Value = "2016-09-05 16:46:41"

# --- Widen the Status / Datetime columns to fit the new text ----------------
# The COM layer quantizes ColumnWidth to 1/6-character steps, so feed it a
# value that lands on the step nearest the authored target width
# (~17.216 chars, i.e. stored width 103/6 = 17.1666...).
$newColWidth = 16.333333333333336

$overview.Columns.Item(5).ColumnWidth = $newColWidth   # column E (zh-cn status)
$overview.Columns.Item(6).ColumnWidth = $newColWidth   # column F (de-de status)
$zhcn.Columns.Item(3).ColumnWidth     = $newColWidth   # column C (Status)
$dede.Columns.Item(3).ColumnWidth     = $newColWidth   # column C (Status)
